$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# 1) "Wechsel von ..." paragraph - collapse split runs / remove proofErr markers
ReplaceText "Wechsel von gefloateden Elementen zu felx-box (grird)" "Wechsel von gefloateden Elementen zu felx-box (grird)"

# 2) "Anpassung der ..." paragraph, first run group (before the inline image)
ReplaceText "Anpassung der reflex.min: auskommentieren von box-sizing, da das den percentage circle zerstört." "Anpassung der reflex.min: auskommentieren von box-sizing, da das den percentage circle zerstört."

# 3) " Dies scheint ein inherit Problem ..." run group (after the inline image)
ReplaceText " Dies scheint ein inherit Problem zu sein – nach einer anderen Lösungsmöglichkeit suchen!!" " Dies scheint ein inherit Problem zu sein – nach einer anderen Lösungsmöglichkeit suchen!!"

# 4) " Folge: refelx.min ..." run group
ReplaceText " Folge: refelx.min ist nicht mehr original. Das muss vermerkt werden." " Folge: refelx.min ist nicht mehr original. Das muss vermerkt werden."

# 5) "Meinestyles.css ..." paragraph (leading run group only, rest stays as-is)
ReplaceText "Meinestyles.css nach Praktikum bei nemo überarbeitungsbedürftig. Wird zukünftig ersetzt durch" "Meinestyles.css nach Praktikum bei nemo überarbeitungsbedürftig. Wird zukünftig ersetzt durch"

# 6) "Aus Kompatibilitätsgründen ..." paragraph
ReplaceText "Aus Kompatibilitätsgründen in Zukunft keine Google Fonts mehr – Helvetica oder Verdana sollten ein ähnliches Ergebnis erzielen – keep it simple!" "Aus Kompatibilitätsgründen in Zukunft keine Google Fonts mehr – Helvetica oder Verdana sollten ein ähnliches Ergebnis erzielen – keep it simple!"

# 7) "Grid Klassen ..." paragraph - merge runs, add trailing period, drop the old
#    _GoBack bookmark (it gets recreated at the very end of the document below)
ReplaceText "Grid Klassen mit Angaben, wie xs usw. versehen um Kompatibilität auf diversen Monitoren zu gewährleisten." "Grid Klassen mit Angaben, wie xs usw. versehen um Kompatibilität auf diversen Monitoren zu gewährleisten."

# 8) Append the three new list items after "Grid Klassen ..."
$gridParaIndex = $d.Paragraphs.Count
$gridPara = $d.Paragraphs($gridParaIndex)
$gridPara.Range.InsertParagraphAfter()
$d.Paragraphs($gridParaIndex + 1).Range.Text = "Reflex grid in integrated frameworks and projects einfügen. Lizenz?!"

$favicon1Index = $gridParaIndex + 1
$favicon1Para = $d.Paragraphs($favicon1Index)
$favicon1Para.Range.InsertParagraphAfter()
$d.Paragraphs($favicon1Index + 1).Range.Text = "Favicon in Wurzelverzeichnis verschoben"

$favicon2Index = $favicon1Index + 1
$favicon2Para = $d.Paragraphs($favicon2Index)
$favicon2Para.Range.InsertParagraphAfter()
# Add a temporary trailing marker character so we can anchor a zero-width
# bookmark exactly at "end of text" without hitting the degenerate
# (Start == End == document end) Range quirk, then strip the marker again.
$d.Paragraphs($favicon2Index + 1).Range.Text = "Favicon ist kaputt. Einfach neu generieren in 32x32#"

$lastPara = $d.Paragraphs($favicon2Index + 1)
$markerPos = $lastPara.Range.End - 2
$markerRange = $d.Range($markerPos, $markerPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$lastPara2 = $d.Paragraphs($favicon2Index + 1)
$deletePos = $lastPara2.Range.End - 1
$deleteRange = $d.Range($deletePos - 1, $deletePos)
$deleteRange.Text = ""
